$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for all data rows (2 through 498)
# from 2023-09-17 (45186) to 2023-09-19 (45188).
$ws.Range("C2:C498").Value = 45188

# Ensure row 498 keeps an explicit row height (matches new trailing rows).
$ws.Rows.Item(498).RowHeight = 15

# Append the new record as row 499.
$ws.Range("A499").Value = "A 43753-2023"
$ws.Range("B499").Value = 45187
$ws.Range("C499").Value = 45188
$ws.Range("D499").Value = "VÄSTMANLANDS LÄN"
$ws.Range("E499").Value = "VÄSTERÅS"
$ws.Range("G499").Value = 1.8
$ws.Range("H499").Value = 0
$ws.Range("I499").Value = 0
$ws.Range("J499").Value = 0
$ws.Range("K499").Value = 0
$ws.Range("L499").Value = 0
$ws.Range("M499").Value = 0
$ws.Range("N499").Value = 0
$ws.Range("O499").Value = 0
$ws.Range("P499").Value = 0
$ws.Range("Q499").Value = 0
$ws.Range("R499").Value = ""

# Match formatting used by the rest of the table (date format for B/C, wrap text for R).
$ws.Range("B499:C499").NumberFormat = $ws.Range("B498:C498").NumberFormat
$ws.Range("R499").WrapText = $true
